$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53. This shifts the existing rows 53-144
# down to 54-145, preserving all of their data/formatting.
$ws.Rows("53:53").Insert()

# Populate the newly inserted row 53 with the new price-record data.
$ws.Cells.Item(53, 1).Value2 = 10
$ws.Cells.Item(53, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(53, 3).Value2 = "La Araucanía"
$ws.Cells.Item(53, 4).Value2 = 44557
$ws.Cells.Item(53, 5).Value2 = 9
$ws.Cells.Item(53, 6).Value2 = "Fruta"
$ws.Cells.Item(53, 7).Value2 = 100103
$ws.Cells.Item(53, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(53, 9).Value2 = 100103002
$ws.Cells.Item(53, 10).Value2 = "Ciruela"
$ws.Cells.Item(53, 11).Value2 = "Red Beaut"
$ws.Cells.Item(53, 12).Value2 = "Primera"
$ws.Cells.Item(53, 13).Value2 = 55
$ws.Cells.Item(53, 14).Value2 = 19000
$ws.Cells.Item(53, 15).Value2 = 20000
$ws.Cells.Item(53, 16).Value2 = 19455
$ws.Cells.Item(53, 17).Value2 = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(53, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(53, 19).Value2 = 1081
$ws.Cells.Item(53, 20).Value2 = 18
